# B1--and-B2-PowerPoint.pptx edit
#
#  1) The table on slide 5 gets a different built-in table style applied
#     (PowerPoint Table Design gallery -> a style whose StyleId is
#     {6AC34C94-4945-4673-8A58-58406140BBEB}, a standard style that is not
#     one of the styles declared in this deck's ppt/tableStyles.xml).
#
#  2) The deck's theme colour palette is switched from the custom
#     "Integral / Red Violet" palette over to the standard Office palette
#     (Design tab -> Variants / Colors -> "Office").

$p = $ppt.ActivePresentation

# --- 1) Re-style the table on slide 5 -------------------------------------
$slide = $p.Slides.Item(5)
$tableShape = $slide.Shapes.Item(2)
$table = $tableShape.Table
$table.ApplyStyle("{6AC34C94-4945-4673-8A58-58406140BBEB}")

# --- 2) Swap the theme colour scheme over to the "Office" palette --------
function ConvertHexToRgb([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $b * 65536 + $g * 256 + $r
}

# Order matches MsoThemeColorSchemeIndex 1..12:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

for ($i = 1; $i -le $officeColors.Length; $i++) {
    $colorScheme.Item($i).RGB = ConvertHexToRgb($officeColors[$i - 1])
}
